$d = $word.ActiveDocument

# The document has a single section whose primary ("default") header/footer
# and "first page" header/footer each carry one inline picture (the BTec
# logo in the headers, the Pearson logo in the footers). Word numbers
# Headers/Footers collections 1 = primary, 2 = first page, 3 = even page.
#
# Rename every one of those inline pictures:
#   - BTec logo pictures (headers):  image1.jpg -> image2.jpg
#   - Pearson logo pictures (footers): image2.png -> image1.png

foreach ($sec in $d.Sections) {
    for ($i = 1; $i -le 3; $i++) {
        $hdr = $sec.Headers.Item($i)
        if ($hdr.Exists) {
            for ($j = 1; $j -le $hdr.Range.InlineShapes.Count; $j++) {
                $shp = $hdr.Range.InlineShapes.Item($j)
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image2.jpg"
                }
            }
        }

        $ftr = $sec.Footers.Item($i)
        if ($ftr.Exists) {
            for ($j = 1; $j -le $ftr.Range.InlineShapes.Count; $j++) {
                $shp = $ftr.Range.InlineShapes.Item($j)
                if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $shp.Name = "image1.png"
                }
            }
        }
    }
}
